$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Anthony M C Alexander"
$ws.Range("B2").Value = "aalexander@dwpv.com"
$ws.Range("C2").Value = "'4163676920"
$ws.Range("D2").Value = "Canada"
$ws.Range("F2").Value = "https://www.dwpv.com/our-people/anthony-alexander"
$ws.Range("I2").Value = "Partner"
$ws.Range("J2").Value = "Davies Ward Phillips And Vineberg"
